# Adds "Progreso Día 2" entries to the tracking sheet:
#  - turns the old blank spacer row (18) into a second header row
#  - fills in five new task rows (19-23) describing Day 2 progress
#  - removes the now-unused merge on the old spacer row
#  - resizes a couple of header rows to fit the new (wrapped) text
#  - relocates the trailing placeholder cell from row 27 to row 28
#  - updates the active selection/scroll position to the new area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn row 18 into the "Progreso Día 2" column header row -----------
# Row 2 already carries the exact header labels/format we need, so copy
# it wholesale (values + formats) onto row 18, then drop the merge that
# used to span A18:E18 as an empty banner cell.
$ws.Range("A2:E2").Copy()
$ws.Range("A18:E18").PasteSpecial(-4104)
$ws.Range("A18:E18").UnMerge()

# --- New data rows 19-23 -------------------------------------------------
# Rows 3-15 already use the per-column formats we need (date column A,
# wrapped/bordered text columns B:E), so clone row 3's formatting across
# each new row before writing the values.
$ws.Range("A3:E3").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)

# Values are written in the same order the author originally typed them
# (row 19, then 21, 22, 23, and finally 20 - a bug note slotted in
# afterwards) so new shared-string entries land at the same indices as
# the canonical file.
$ws.Range("A19").Value = 45493
$ws.Range("B19").Value = "Corrección de filtro de búsqueda"
$ws.Range("C19").Value = "Mejora en la funcionalidad del filtro de búsqueda para incluir la capital del país en el criterio de búsqueda."
$ws.Range("D19").Value = "src/components/CountryList.jsx"
$ws.Range("E19").Value = "La búsqueda ahora considera tanto el nombre como la capital del país."

$ws.Range("A21").Value = 45493
$ws.Range("B21").Value = "Añadir botón de retroceso en CountryList"
$ws.Range("C21").Value = "Añadido un botón para regresar a la lista de países desde la vista de detalles de país."
$ws.Range("D21").Value = "src/components/CountryList.jsx"
$ws.Range("E21").Value = "Facilita la navegación desde el detalle del país a la lista de países."

$ws.Range("A22").Value = 45493
$ws.Range("B22").Value = "Añadir botón de retroceso en CountryDetail"
$ws.Range("C22").Value = "Añadido un botón para regresar a la lista de países desde la vista de detalles de país."
$ws.Range("D22").Value = "src/components/CountryDetail.jsx"
$ws.Range("E22").Value = "Facilita la navegación desde el detalle del país a la lista de países."

$ws.Range("A23").Value = 45493
$ws.Range("B23").Value = "Añadir botón de retroceso en NotFound"
$ws.Range("C23").Value = "Añadido un botón para regresar a la página de inicio desde la vista de ""Página No Encontrada""."
$ws.Range("D23").Value = "src/components/NotFound.jsx"
$ws.Range("E23").Value = "Proporciona una forma de regresar a la página principal desde ""No Encontrado""."

$ws.Range("A20").Value = 45493
$ws.Range("B20").Value = "Corrección de Errores"
$ws.Range("C20").Value = "Se encontró error grave, no existía el CRUD que se creía agregado en el día uno, no se podía agregar, eliminar o actualizar países en la LocalStorage, solo se puede ver los países de la API, se intentó solucionar agregándolo, pero hubieron muchos más errores, por el momento no se logró solucionar."
$ws.Range("D20").Value = "src/components/CountryList.jsx"
$ws.Range("E20").Value = "No se logró implementar el CRUD completo, se identificaron problemas adicionales."

# --- Move the trailing placeholder cell from row 27 down to row 28 -----
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A27").Clear()

# --- Row heights ----------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 42
$ws.Rows.Item(18).RowHeight = 39
$ws.Rows.Item(19).RowHeight = 105
$ws.Rows.Item(20).RowHeight = 180
$ws.Rows.Item(21).RowHeight = 90
$ws.Rows.Item(22).RowHeight = 114
$ws.Rows.Item(23).RowHeight = 105

# --- Selection / scroll position matching the saved view ----------------
$ws.Range("G23").Select()
$excel.ActiveWindow.ScrollRow = 21
